# Auto-generated script applying scheduled market-data refresh to the leve profit sheets.
# For each sheet, columns H-N (currentAveragePrice.. LeveProfitHQ) are refreshed with new values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 155731.53
$ws.Range("I86").Value = 1914.7142
$ws.Range("J86").Value = 335184.5
$ws.Range("K86").Value = 1914.7142
$ws.Range("L86").Value = 335184.5
$ws.Range("M86").Value = -791.7141999999999
$ws.Range("N86").Value = -337430.5
$ws.Range("H89").Value = 155731.53
$ws.Range("I89").Value = 1914.7142
$ws.Range("J89").Value = 335184.5
$ws.Range("K89").Value = 9573.571
$ws.Range("L89").Value = 1675922.5
$ws.Range("M89").Value = -3957.571
$ws.Range("N89").Value = -1687154.5
$ws.Range("H101").Value = 20003978
$ws.Range("J101").Value = 999
$ws.Range("L101").Value = 2997
$ws.Range("N101").Value = -6241
$ws.Range("H111").Value = 30531.834
$ws.Range("I111").Value = 28266.445
$ws.Range("J111").Value = 37328
$ws.Range("K111").Value = 84799.33499999999
$ws.Range("L111").Value = 111984
$ws.Range("M111").Value = -81732.33499999999
$ws.Range("N111").Value = -118118
$ws.Range("H113").Value = 4250
$ws.Range("J113").Value = 4500
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -11008
$ws.Range("H132").Value = 2448.111
$ws.Range("I132").Value = 2448.111
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7344.333
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4814.333
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 15000
$ws.Range("I55").Value = 15000
$ws.Range("K55").Value = 15000
$ws.Range("M55").Value = -14685
$ws.Range("H61").Value = 1498
$ws.Range("I61").Value = 1498
$ws.Range("K61").Value = 1498
$ws.Range("M61").Value = -1286
$ws.Range("H110").Value = 4626465
$ws.Range("I110").Value = 5286960
$ws.Range("K110").Value = 5286960
$ws.Range("M110").Value = -5284915
$ws.Range("H122").Value = 1769.6666
$ws.Range("I122").Value = 1769.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5308.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2858.9998
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 1498
$ws.Range("I136").Value = 1498
$ws.Range("K136").Value = 4494
$ws.Range("M136").Value = -1944

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3204
$ws.Range("I107").Value = 4810.6665
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 4810.6665
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = -2890.6665
$ws.Range("N107").Value = -5839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1903.4286
$ws.Range("I58").Value = 1665.1666
$ws.Range("J58").Value = 3333
$ws.Range("K58").Value = 1665.1666
$ws.Range("L58").Value = 3333
$ws.Range("M58").Value = -1462.1666
$ws.Range("N58").Value = -3739
$ws.Range("H62").Value = 2800
$ws.Range("I62").Value = 2800
$ws.Range("K62").Value = 2800
$ws.Range("M62").Value = -2176
$ws.Range("H65").Value = 2800
$ws.Range("I65").Value = 2800
$ws.Range("K65").Value = 14000
$ws.Range("M65").Value = -10880
$ws.Range("H75").Value = 34999
$ws.Range("J75").Value = 34999
$ws.Range("L75").Value = 34999
$ws.Range("N75").Value = -36995
$ws.Range("H78").Value = 34999
$ws.Range("J78").Value = 34999
$ws.Range("L78").Value = 104997
$ws.Range("N78").Value = -114981
$ws.Range("H99").Value = 1984.5385
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 2449.5
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 2449.5
$ws.Range("M99").Value = -402
$ws.Range("N99").Value = -5445.5
$ws.Range("H122").Value = 2854.3635
$ws.Range("I122").Value = 1483
$ws.Range("K122").Value = 4449
$ws.Range("M122").Value = -1999
$ws.Range("H126").Value = 1984.5385
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 2449.5
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 7348.5
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -12288.5
$ws.Range("H136").Value = 1903.4286
$ws.Range("I136").Value = 1665.1666
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 4995.4998
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -2445.4998
$ws.Range("N136").Value = -15099

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3600
$ws.Range("J75").Value = 3133.3333
$ws.Range("L75").Value = 9399.999899999999
$ws.Range("N75").Value = -11395.9999
$ws.Range("H78").Value = 3600
$ws.Range("J78").Value = 3133.3333
$ws.Range("L78").Value = 28199.9997
$ws.Range("N78").Value = -38183.9997
$ws.Range("H119").Value = 1404.8
$ws.Range("I119").Value = 1404.8
$ws.Range("K119").Value = 4214.4
$ws.Range("M119").Value = 623.6000000000004
$ws.Range("H131").Value = 437003.75
$ws.Range("I131").Value = 1427.5
$ws.Range("J131").Value = 590736.5
$ws.Range("K131").Value = 4282.5
$ws.Range("L131").Value = 1772209.5
$ws.Range("M131").Value = 757.5
$ws.Range("N131").Value = -1782289.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 25772.285
$ws.Range("J94").Value = 27962.818
$ws.Range("L94").Value = 27962.818
$ws.Range("N94").Value = -29314.818
$ws.Range("H102").Value = 1921.3846
$ws.Range("I102").Value = 1897.8
$ws.Range("K102").Value = 1897.8
$ws.Range("M102").Value = -275.8
$ws.Range("H107").Value = 1538.0667
$ws.Range("I107").Value = 754.7778
$ws.Range("J107").Value = 2713
$ws.Range("K107").Value = 754.7778
$ws.Range("L107").Value = 2713
$ws.Range("M107").Value = 1165.2222
$ws.Range("N107").Value = -6553
$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 999
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5339
$ws.Range("H122").Value = 5934
$ws.Range("I122").Value = 5256.6665
$ws.Range("K122").Value = 15769.9995
$ws.Range("M122").Value = -13319.9995
$ws.Range("H126").Value = 2389
$ws.Range("I126").Value = 1986.25
$ws.Range("K126").Value = 5958.75
$ws.Range("M126").Value = -3488.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7490.1763
$ws.Range("I7").Value = 6512.7144
$ws.Range("J7").Value = 8174.4
$ws.Range("K7").Value = 6512.7144
$ws.Range("L7").Value = 8174.4
$ws.Range("M7").Value = -6400.7144
$ws.Range("N7").Value = -8398.4
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 809
$ws.Range("I61").Value = 799
$ws.Range("J61").Value = 824
$ws.Range("K61").Value = 799
$ws.Range("L61").Value = 824
$ws.Range("M61").Value = -597
$ws.Range("N61").Value = -1228
$ws.Range("H82").Value = 509.55554
$ws.Range("I82").Value = 378.6
$ws.Range("J82").Value = 673.25
$ws.Range("K82").Value = 378.6
$ws.Range("L82").Value = 673.25
$ws.Range("M82").Value = -17.60000000000002
$ws.Range("N82").Value = -1395.25
$ws.Range("H85").Value = 509.55554
$ws.Range("I85").Value = 378.6
$ws.Range("J85").Value = 673.25
$ws.Range("K85").Value = 378.6
$ws.Range("L85").Value = 673.25
$ws.Range("M85").Value = 869.4
$ws.Range("N85").Value = -3169.25
$ws.Range("H113").Value = 809
$ws.Range("I113").Value = 799
$ws.Range("J113").Value = 824
$ws.Range("K113").Value = 799
$ws.Range("L113").Value = 824
$ws.Range("M113").Value = 1371
$ws.Range("N113").Value = -5164
$ws.Range("H126").Value = 7490.1763
$ws.Range("I126").Value = 6512.7144
$ws.Range("J126").Value = 8174.4
$ws.Range("K126").Value = 19538.1432
$ws.Range("L126").Value = 24523.2
$ws.Range("M126").Value = -17068.1432
$ws.Range("N126").Value = -29463.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 29750
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 38000
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 38000
$ws.Range("M82").Value = -4617
$ws.Range("N82").Value = -38766
$ws.Range("H85").Value = 29750
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 38000
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 38000
$ws.Range("M85").Value = -3674
$ws.Range("N85").Value = -40652
$ws.Range("H126").Value = 3381.25
$ws.Range("I126").Value = 1963.25
$ws.Range("J126").Value = 4799.25
$ws.Range("K126").Value = 5889.75
$ws.Range("L126").Value = 14397.75
$ws.Range("M126").Value = -3419.75
$ws.Range("N126").Value = -19337.75
